$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 9f9a1b7b...md row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-11-29 04:13:43"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" for 9f9a1b7b...md row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-11-29 04:13:27"
$wsZhCn.Range("K3").Value = "2016-11-29 04:14:23"

# de-de sheet: "Correspond Handoff Datetime" (shares same underlying value as
# Overview!G3 above) and "Correspond Handback DateTime" for 9f9a1b7b...md row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-11-29 04:13:43"
$wsDeDe.Range("K3").Value = "2016-11-29 04:14:41"
